# Auto-generated edit script applying numeric corrections to the
# Leve profit-tracking workbook (columns H-N per sheet table).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3846.4285
$ws.Range("J17").Value = 3846.4285
$ws.Range("L17").Value = 11539.2855
$ws.Range("N17").Value = -11875.2855
$ws.Range("H33").Value = 505.8889
$ws.Range("I33").Value = 360.83334
$ws.Range("K33").Value = 360.83334
$ws.Range("M33").Value = -131.83334
$ws.Range("H44").Value = 29500
$ws.Range("I44").Value = 24000
$ws.Range("J44").Value = 35000
$ws.Range("K44").Value = 24000
$ws.Range("L44").Value = 35000
$ws.Range("M44").Value = -23538
$ws.Range("N44").Value = -35924
$ws.Range("H58").Value = 3875.5715
$ws.Range("I58").Value = 15
$ws.Range("J58").Value = 4519
$ws.Range("K58").Value = 45
$ws.Range("L58").Value = 13557
$ws.Range("M58").Value = 105
$ws.Range("N58").Value = -13857
$ws.Range("H80").Value = 2714.5952
$ws.Range("I80").Value = 1314.5834
$ws.Range("J80").Value = 3274.6
$ws.Range("K80").Value = 3943.7502
$ws.Range("L80").Value = 9823.799999999999
$ws.Range("M80").Value = -2945.7502
$ws.Range("N80").Value = -11819.8
$ws.Range("H83").Value = 2714.5952
$ws.Range("I83").Value = 1314.5834
$ws.Range("J83").Value = 3274.6
$ws.Range("K83").Value = 11831.2506
$ws.Range("L83").Value = 29471.4
$ws.Range("M83").Value = -6839.250599999999
$ws.Range("N83").Value = -39455.39999999999
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22496
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -72480
$ws.Range("H101").Value = 2621
$ws.Range("I101").Value = 992.3333
$ws.Range("J101").Value = 4249.6665
$ws.Range("K101").Value = 2976.9999
$ws.Range("L101").Value = 12748.9995
$ws.Range("M101").Value = -1354.9999
$ws.Range("N101").Value = -15992.9995
$ws.Range("H137").Value = 1892.5555
$ws.Range("I137").Value = 1414.8125
$ws.Range("J137").Value = 2587.4546
$ws.Range("K137").Value = 4244.4375
$ws.Range("L137").Value = 7762.3638
$ws.Range("M137").Value = -1694.4375
$ws.Range("N137").Value = -12862.3638

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 1050
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H32").Value = 16488.197
$ws.Range("I32").Value = 12348.726
$ws.Range("K32").Value = 12348.726
$ws.Range("M32").Value = -12061.726
$ws.Range("H110").Value = 3442.2632
$ws.Range("I110").Value = 3442.2632
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 3442.2632
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -1397.2632
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 3343.8
$ws.Range("I122").Value = 2462.7827
$ws.Range("J122").Value = 5032.4165
$ws.Range("K122").Value = 7388.348100000001
$ws.Range("L122").Value = 15097.2495
$ws.Range("M122").Value = -4938.348100000001
$ws.Range("N122").Value = -19997.2495
$ws.Range("H130").Value = 292286
$ws.Range("J130").Value = 292286
$ws.Range("L130").Value = 292286
$ws.Range("N130").Value = -302326
$ws.Range("H132").Value = 3930.0417
$ws.Range("I132").Value = 3661.2896
$ws.Range("K132").Value = 10983.8688
$ws.Range("M132").Value = -8453.8688

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3414.7693
$ws.Range("I105").Value = 2323.3809
$ws.Range("J105").Value = 7998.6
$ws.Range("K105").Value = 2323.3809
$ws.Range("L105").Value = 7998.6
$ws.Range("M105").Value = -576.3809000000001
$ws.Range("N105").Value = -11492.6
$ws.Range("H107").Value = 545.375
$ws.Range("I107").Value = 527
$ws.Range("K107").Value = 527
$ws.Range("M107").Value = 1393
$ws.Range("H134").Value = 2725.9546
$ws.Range("I134").Value = 1726.2222
$ws.Range("K134").Value = 5178.6666
$ws.Range("M134").Value = -2643.6666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 485.5
$ws.Range("I16").Value = 485.5
$ws.Range("K16").Value = 485.5
$ws.Range("M16").Value = -198.5
$ws.Range("H31").Value = 3242.1792
$ws.Range("J31").Value = 4338.7417
$ws.Range("L31").Value = 4338.7417
$ws.Range("N31").Value = -4928.7417
$ws.Range("H34").Value = 3242.1792
$ws.Range("J34").Value = 4338.7417
$ws.Range("L34").Value = 4338.7417
$ws.Range("N34").Value = -4742.7417
$ws.Range("H53").Value = 75000
$ws.Range("J53").Value = 75000
$ws.Range("L53").Value = 75000
$ws.Range("N53").Value = -76214
$ws.Range("H57").Value = 42500
$ws.Range("I57").Value = 35000
$ws.Range("J57").Value = 50000
$ws.Range("K57").Value = 35000
$ws.Range("L57").Value = 50000
$ws.Range("M57").Value = -34440
$ws.Range("N57").Value = -51120
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H113").Value = 485.5
$ws.Range("I113").Value = 485.5
$ws.Range("K113").Value = 485.5
$ws.Range("M113").Value = 1684.5
$ws.Range("H132").Value = 2270.9375
$ws.Range("I132").Value = 1913.0344
$ws.Range("J132").Value = 5730.6665
$ws.Range("K132").Value = 5739.1032
$ws.Range("L132").Value = 17191.9995
$ws.Range("M132").Value = -3209.1032
$ws.Range("N132").Value = -22251.9995
$ws.Range("H134").Value = 2850.4243
$ws.Range("I134").Value = 1757.8572
$ws.Range("J134").Value = 4762.4165
$ws.Range("K134").Value = 5273.571599999999
$ws.Range("L134").Value = 14287.2495
$ws.Range("M134").Value = -2738.571599999999
$ws.Range("N134").Value = -19357.2495

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 248.8
$ws.Range("I2").Value = 137.88889
$ws.Range("J2").Value = 415.16666
$ws.Range("K2").Value = 827.33334
$ws.Range("L2").Value = 2490.99996
$ws.Range("M2").Value = -714.33334
$ws.Range("N2").Value = -2716.99996
$ws.Range("H38").Value = 825.71875
$ws.Range("J38").Value = 1560.8125
$ws.Range("L38").Value = 4682.4375
$ws.Range("N38").Value = -5376.4375
$ws.Range("H50").Value = 970.7143
$ws.Range("I50").Value = 132.5
$ws.Range("J50").Value = 6000
$ws.Range("K50").Value = 397.5
$ws.Range("L50").Value = 18000
$ws.Range("M50").Value = 83.5
$ws.Range("N50").Value = -18962
$ws.Range("H53").Value = 970.7143
$ws.Range("I53").Value = 132.5
$ws.Range("J53").Value = 6000
$ws.Range("K53").Value = 397.5
$ws.Range("L53").Value = 18000
$ws.Range("M53").Value = 83.5
$ws.Range("N53").Value = -18962
$ws.Range("H55").Value = 27782678
$ws.Range("J55").Value = 27782678
$ws.Range("L55").Value = 83348034
$ws.Range("N55").Value = -83348388
$ws.Range("H70").Value = 77017.664
$ws.Range("I70").Value = 2531.4
$ws.Range("K70").Value = 7594.200000000001
$ws.Range("M70").Value = -7279.200000000001
$ws.Range("H73").Value = 77017.664
$ws.Range("I73").Value = 2531.4
$ws.Range("K73").Value = 7594.200000000001
$ws.Range("M73").Value = -6502.200000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 6395537
$ws.Range("J20").Value = 50129.715
$ws.Range("L20").Value = 50129.715
$ws.Range("N20").Value = -50619.715
$ws.Range("H132").Value = 4947.622
$ws.Range("I132").Value = 4126.0835
$ws.Range("J132").Value = 5886.524
$ws.Range("K132").Value = 12378.2505
$ws.Range("L132").Value = 17659.572
$ws.Range("M132").Value = -9848.250499999998
$ws.Range("N132").Value = -22719.572

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 2027242.8
$ws.Range("I43").Value = 38000
$ws.Range("J43").Value = 2524553.5
$ws.Range("K43").Value = 38000
$ws.Range("L43").Value = 2524553.5
$ws.Range("M43").Value = -37807
$ws.Range("N43").Value = -2524939.5
$ws.Range("H46").Value = 9291.638999999999
$ws.Range("I46").Value = 4916.5
$ws.Range("K46").Value = 4916.5
$ws.Range("M46").Value = -4728.5
$ws.Range("H61").Value = 2512.4546
$ws.Range("I61").Value = 1496.6923
$ws.Range("J61").Value = 6285.2856
$ws.Range("K61").Value = 1496.6923
$ws.Range("L61").Value = 6285.2856
$ws.Range("M61").Value = -1294.6923
$ws.Range("N61").Value = -6689.2856
$ws.Range("H113").Value = 2512.4546
$ws.Range("I113").Value = 1496.6923
$ws.Range("J113").Value = 6285.2856
$ws.Range("K113").Value = 1496.6923
$ws.Range("L113").Value = 6285.2856
$ws.Range("M113").Value = 673.3077000000001
$ws.Range("N113").Value = -10625.2856
$ws.Range("H132").Value = 4690.2354
$ws.Range("I132").Value = 4024.4075
$ws.Range("J132").Value = 7258.4287
$ws.Range("K132").Value = 12073.2225
$ws.Range("L132").Value = 21775.2861
$ws.Range("M132").Value = -9543.2225
$ws.Range("N132").Value = -26835.2861

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 9900
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 9900
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 9900
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -10124
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H96").Value = 102079
$ws.Range("I96").Value = 250951.5
$ws.Range("J96").Value = 2830.6667
$ws.Range("K96").Value = 250951.5
$ws.Range("L96").Value = 2830.6667
$ws.Range("M96").Value = -249578.5
$ws.Range("N96").Value = -5576.6667
$ws.Range("H100").Value = 797.2857
$ws.Range("I100").Value = 763.4
$ws.Range("J100").Value = 882
$ws.Range("K100").Value = 1526.8
$ws.Range("L100").Value = 1764
$ws.Range("M100").Value = -985.8
$ws.Range("N100").Value = -2846
$ws.Range("H132").Value = 2804.5757
$ws.Range("I132").Value = 1430.8
$ws.Range("K132").Value = 4292.4
$ws.Range("M132").Value = -1762.4

